$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.456.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.77%  "

$ws.Range("D3").Value = "'1.858.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.67%  "

$ws.Range("E4").Value = "  -0.25%  "

$ws.Range("D5").Value = "'311.54"
$ws.Range("D5").Style = "Normal"

$ws.Range("E6").Value = "  -0.22%  "

$ws.Range("D7").Value = "'0.4763"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "'0.3793"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.07%  "

$ws.Range("D9").Value = "'0.07311"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.17%  "

$ws.Range("D10").Value = "'0.9296"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.23%  "

$ws.Range("D12").Value = "'0.07789"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.75%  "

$ws.Range("D13").Value = "'1.852.97"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.27%  "

$ws.Range("D14").Value = "'5.450"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.33%  "

$ws.Range("E15").Value = "  +1.38%  "

$ws.Range("D16").Value = "'90.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.37%  "

$ws.Range("E17").Value = "  -0.37%  "

$ws.Range("D18").Value = "'0.000008814"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.72%  "

$ws.Range("E19").Value = "  -0.35%  "

$ws.Range("D20").Value = "'27.479.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.70%  "

$ws.Range("D21").Value = "'14.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.38%  "

$ws.Range("D22").Value = "'5.091"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.39%  "

$ws.Range("E23").Value = "  +0.52%  "

$ws.Range("D24").Value = "'1.939"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.99%  "

$ws.Range("D25").Value = "'154.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.27%  "

$ws.Range("E26").Value = "  +1.11%  "

$ws.Range("D27").Value = "'2.006"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.28%  "

$ws.Range("E28").Value = "  +0.71%  "

$ws.Range("D29").Value = "'4.941"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.65%  "

$ws.Range("D30").Value = "'0.08875"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.09%  "

$ws.Range("D31").Value = "'3.329"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.67%  "

$ws.Range("E32").Value = "  +1.97%  "

$ws.Range("D33").Value = "'0.7505"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.87%  "

$ws.Range("D34").Value = "'4.574"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.53%  "

$ws.Range("D35").Value = "'2.699"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.35%  "

$ws.Range("D36").Value = "'0.02045"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.15%  "

$ws.Range("D37").Value = "'1.121"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.51%  "

$ws.Range("D38").Value = "'0.5559"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.72%  "

$ws.Range("D39").Value = "'0.05279"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.17%  "

$ws.Range("D40").Value = "'2.981"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.42%  "

$ws.Range("D41").Value = "'7.004"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.27%  "

$ws.Range("D42").Value = "'8.555"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.11%  "

$ws.Range("E44").Value = "  +2.65%  "

$ws.Range("E45").Value = "  +0.29%  "

$ws.Range("E46").Value = "  -0.33%  "

$ws.Range("D47").Value = "'103.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.69%  "

$ws.Range("D48").Value = "'1.663"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.33%  "

$ws.Range("D49").Value = "'67.22"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.98%  "

$ws.Range("D50").Value = "'0.06090"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.20%  "

$ws.Range("D51").Value = "'0.9111"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.13%  "
